$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (12) so that the existing
# "Notes on data quality" / "Other notes" columns shift right by one
# (old L -> M, old M -> N).
$ws.Columns.Item(12).Insert()

# New header cell (L1) - "audience" question. Copy formatting from the
# neighboring header cell (K1) first, then set the text.
$ws.Cells.Item(1, 11).Copy()
$ws.Cells.Item(1, 12).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 12).Value = "Who is the audience for this data?"

# New data cells for rows 2 and 3 (L2, L3), matching the style used by the
# rest of that row (copy format from K2 / K3).
$ws.Cells.Item(2, 11).Copy()
$ws.Cells.Item(2, 12).PasteSpecial(-4122)
$ws.Cells.Item(2, 12).Value = "researchers, policy makers"

$ws.Cells.Item(3, 11).Copy()
$ws.Cells.Item(3, 12).PasteSpecial(-4122)
$ws.Cells.Item(3, 12).Value = "researchers, policy makers"

$ws.Application.CutCopyMode = 0

# Adjust row heights for the rows whose content grew.
$ws.Rows.Item(2).RowHeight = 102
$ws.Rows.Item(3).RowHeight = 301

# Update the view: select the newly added last cell (N3) and scroll down
# so row 3 is at the top of the viewport.
$ws.Activate()
$ws.Range("N3").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
